$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 522.25806
$ws.Range("I28").Value = 336.9565
$ws.Range("J28").Value = 1055
$ws.Range("K28").Value = 336.9565
$ws.Range("L28").Value = 1055
$ws.Range("M28").Value = 148.0435
$ws.Range("N28").Value = -2025
$ws.Range("H62").Value = 14289.5
$ws.Range("I62").Value = 4447.5
$ws.Range("J62").Value = 16750
$ws.Range("K62").Value = 4447.5
$ws.Range("L62").Value = 16750
$ws.Range("M62").Value = -3823.5
$ws.Range("N62").Value = -17998
$ws.Range("H65").Value = 14289.5
$ws.Range("I65").Value = 4447.5
$ws.Range("J65").Value = 16750
$ws.Range("K65").Value = 22237.5
$ws.Range("L65").Value = 83750
$ws.Range("M65").Value = -19117.5
$ws.Range("N65").Value = -89990
$ws.Range("H107").Value = 401.72726
$ws.Range("I107").Value = 414.05264
$ws.Range("J107").Value = 323.66666
$ws.Range("K107").Value = 414.05264
$ws.Range("L107").Value = 323.66666
$ws.Range("M107").Value = 1505.94736
$ws.Range("N107").Value = -4163.66666
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("H137").Value = 5800.6855
$ws.Range("I137").Value = 5460.2593
$ws.Range("J137").Value = 6949.625
$ws.Range("K137").Value = 16380.7779
$ws.Range("L137").Value = 20848.875
$ws.Range("M137").Value = -13830.7779
$ws.Range("N137").Value = -25948.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1009.4737
$ws.Range("I2").Value = 904.8570999999999
$ws.Range("J2").Value = 1302.4
$ws.Range("K2").Value = 904.8570999999999
$ws.Range("L2").Value = 1302.4
$ws.Range("M2").Value = -791.8570999999999
$ws.Range("N2").Value = -1528.4
$ws.Range("H7").Value = 85198
$ws.Range("J7").Value = 88663.336
$ws.Range("L7").Value = 88663.336
$ws.Range("N7").Value = -88891.336
$ws.Range("H45").Value = 23812682
$ws.Range("I45").Value = 45456630
$ws.Range("K45").Value = 45456630
$ws.Range("M45").Value = -45456253
$ws.Range("H61").Value = 44125580
$ws.Range("I61").Value = 100007440
$ws.Range("J61").Value = 20841474
$ws.Range("K61").Value = 100007440
$ws.Range("L61").Value = 20841474
$ws.Range("M61").Value = -100007228
$ws.Range("N61").Value = -20841898
$ws.Range("H102").Value = 13616
$ws.Range("I102").Value = 13292
$ws.Range("K102").Value = 13292
$ws.Range("M102").Value = -11670
$ws.Range("H110").Value = 4670
$ws.Range("I110").Value = 5000
$ws.Range("J110").Value = 4505
$ws.Range("K110").Value = 5000
$ws.Range("L110").Value = 4505
$ws.Range("M110").Value = -2955
$ws.Range("N110").Value = -8595
$ws.Range("H116").Value = 1009.4737
$ws.Range("I116").Value = 904.8570999999999
$ws.Range("J116").Value = 1302.4
$ws.Range("K116").Value = 904.8570999999999
$ws.Range("L116").Value = 1302.4
$ws.Range("M116").Value = 1389.1429
$ws.Range("N116").Value = -5890.4
$ws.Range("H119").Value = 58333.332
$ws.Range("J119").Value = 58333.332
$ws.Range("L119").Value = 58333.332
$ws.Range("N119").Value = -68009.33199999999
$ws.Range("H136").Value = 44125580
$ws.Range("I136").Value = 100007440
$ws.Range("J136").Value = 20841474
$ws.Range("K136").Value = 300022320
$ws.Range("L136").Value = 62524422
$ws.Range("M136").Value = -300019770
$ws.Range("N136").Value = -62529522

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1009.4737
$ws.Range("I3").Value = 904.8570999999999
$ws.Range("J3").Value = 1302.4
$ws.Range("K3").Value = 904.8570999999999
$ws.Range("L3").Value = 1302.4
$ws.Range("M3").Value = -790.8570999999999
$ws.Range("N3").Value = -1530.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1298.5714
$ws.Range("I16").Value = 889
$ws.Range("K16").Value = 889
$ws.Range("M16").Value = -602
$ws.Range("H22").Value = 412.85715
$ws.Range("I22").Value = 397.6
$ws.Range("K22").Value = 397.6
$ws.Range("M22").Value = -47.60000000000002
$ws.Range("H81").Value = 119000
$ws.Range("J81").Value = 119000
$ws.Range("L81").Value = 119000
$ws.Range("N81").Value = -120996
$ws.Range("H84").Value = 119000
$ws.Range("J84").Value = 119000
$ws.Range("L84").Value = 357000
$ws.Range("N84").Value = -366984
$ws.Range("H107").Value = 3072.5833
$ws.Range("I107").Value = 1539.6
$ws.Range("K107").Value = 1539.6
$ws.Range("M107").Value = 380.4000000000001
$ws.Range("H113").Value = 1298.5714
$ws.Range("I113").Value = 889
$ws.Range("K113").Value = 889
$ws.Range("M113").Value = 1281
$ws.Range("H132").Value = 4429.75
$ws.Range("I132").Value = 2761.077
$ws.Range("K132").Value = 8283.231
$ws.Range("M132").Value = -5753.231
$ws.Range("H134").Value = 3220.5557
$ws.Range("I134").Value = 3245.4119
$ws.Range("K134").Value = 9736.235700000001
$ws.Range("M134").Value = -7201.235700000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 11905057
$ws.Range("I97").Value = 11905057
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 35715171
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -35714675
$ws.Range("N97").ClearContents()
$ws.Range("H109").Value = 1979.5555
$ws.Range("I109").Value = 1545.2858
$ws.Range("K109").Value = 4635.857400000001
$ws.Range("M109").Value = -3595.857400000001
$ws.Range("H132").Value = 1650
$ws.Range("J132").Value = 1726.3636
$ws.Range("L132").Value = 15537.2724
$ws.Range("N132").Value = -20597.2724
$ws.Range("H134").Value = 10259.485
$ws.Range("I134").Value = 1766.125
$ws.Range("K134").Value = 5298.375
$ws.Range("M134").Value = -228.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 96.90909000000001
$ws.Range("J2").Value = 176.8
$ws.Range("L2").Value = 176.8
$ws.Range("N2").Value = -402.8
$ws.Range("H46").Value = 15099.667
$ws.Range("J46").Value = 21499.5
$ws.Range("L46").Value = 21499.5
$ws.Range("N46").Value = -21811.5
$ws.Range("H57").Value = 19999
$ws.Range("J57").Value = 19999
$ws.Range("L57").Value = 19999
$ws.Range("N57").Value = -21639
$ws.Range("H80").Value = 2838.6538
$ws.Range("I80").Value = 2345.2856
$ws.Range("J80").Value = 4910.8
$ws.Range("K80").Value = 2345.2856
$ws.Range("L80").Value = 4910.8
$ws.Range("M80").Value = -1347.2856
$ws.Range("N80").Value = -6906.8
$ws.Range("H83").Value = 2838.6538
$ws.Range("I83").Value = 2345.2856
$ws.Range("J83").Value = 4910.8
$ws.Range("K83").Value = 11726.428
$ws.Range("L83").Value = 24554
$ws.Range("M83").Value = -6734.428
$ws.Range("N83").Value = -34538
$ws.Range("H113").Value = 5001.857
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 62511184
$ws.Range("J132").Value = 22735.166
$ws.Range("L132").Value = 68205.49800000001
$ws.Range("N132").Value = -73265.49800000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2358.5454
$ws.Range("I61").Value = 2438.5557
$ws.Range("K61").Value = 2438.5557
$ws.Range("M61").Value = -2236.5557
$ws.Range("H113").Value = 2358.5454
$ws.Range("I113").Value = 2438.5557
$ws.Range("K113").Value = 2438.5557
$ws.Range("M113").Value = -268.5556999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 28990
$ws.Range("J82").Value = 28990
$ws.Range("L82").Value = 28990
$ws.Range("N82").Value = -29756
$ws.Range("H85").Value = 28990
$ws.Range("J85").Value = 28990
$ws.Range("L85").Value = 28990
$ws.Range("N85").Value = -31642
$ws.Range("H88").Value = 24879
$ws.Range("J88").Value = 24879
$ws.Range("L88").Value = 24879
$ws.Range("N88").Value = -25691
$ws.Range("H91").Value = 24879
$ws.Range("J91").Value = 24879
$ws.Range("L91").Value = 24879
$ws.Range("N91").Value = -27687
$ws.Range("H135").Value = 184357.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 184357.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 184357.25
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -194497.25
$ws.Range("H138").Value = 59749.75
$ws.Range("I138").Value = 58000
$ws.Range("J138").Value = 60333
$ws.Range("K138").Value = 58000
$ws.Range("L138").Value = 60333
$ws.Range("M138").Value = -52860
$ws.Range("N138").Value = -70613
